# adding contributions from dual readout
#
# This script reproduces the authored changes to ListOfGroups.xlsx:
#  - Adds a new "DualReadout" contribution row (row 27) to the "Calice" sheet
#  - Renames the "Calice" sheet to "Calorimetry"
#  - Updates the selections/active sheet to reflect the editing session
#    (US Groups row 6 selected, TPC cell H22 selected, Calorimetry cell H28
#    selected last so it ends up as the active/selected tab)

$wb = $excel.ActiveWorkbook

# --- US Groups sheet: row 6 gets selected at some point during the session ---
$wsUS = $wb.Worksheets.Item("US Groups")
$wsUS.Rows("6:6").Select()

# --- TPC sheet: cursor left on H22 ---
$wsTPC = $wb.Worksheets.Item("TPC")
$wsTPC.Range("H22").Select()

# --- Calice sheet: add the new DualReadout / Iowa State contribution row ---
$wsCalice = $wb.Worksheets.Item("Calice")

$wsCalice.Range("A27").Value = "DualReadout "
$wsCalice.Range("B27").Value = "Iowa State"
$wsCalice.Range("C27").Value = "John Hauptman"
$wsCalice.Range("D27").Value = "hauptman@fnal.gov"

$wsCalice.Range("E27").Value = 41757
$wsCalice.Range("E27").NumberFormat = "d-mmm"

$wsCalice.Range("F27").Value = "Jan"

$wsCalice.Range("G27").Value = 41768
$wsCalice.Range("G27").NumberFormat = "d-mmm"

$wsCalice.Range("H27").Value = "DualReadout_JohnHauptman_IowaState_20140509.pdf"

# Rename the sheet to reflect its broader scope now that Dual Readout is included
$wsCalice.Name = "Calorimetry"

# Leave the selection on the newly added row, which also makes this the
# active/selected sheet tab (last sheet selected wins)
$wsCalice.Range("H28").Select()
